$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_11.1")

# Row 21 used to carry the "Fuente: SICT..." note in C21 (and C22 held the
# "Ultima actualización"/"Dirección General" footer lines). Give B21 the
# same look as B20 (the other footer-label cell) before moving text into it.
$ws.Range("B20").Copy()
$ws.Range("B21").PasteSpecial(-4122)

# Row 20: source note text changes from the old "Fuente..." text
# to the new "Actualización: mayo 2024." text.
$ws.Range("B20").Value = "Actualización: mayo 2024."

# Row 21: B21 now holds the "Fuente: SICT..." text (moved out of B20),
# and C21 loses its old text value (becomes blank, keeps style).
$ws.Range("B21").Value = "Fuente: SICT. Subsecretaria de Infraestructura."
$ws.Range("C21").Value = ""

# Row 22: C22 loses its old text value (becomes blank, keeps style).
$ws.Range("C22").Value = ""
